$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3127.8
$ws.Range("I62").Value = 3000.6667
$ws.Range("K62").Value = 3000.6667
$ws.Range("M62").Value = -2376.6667

$ws.Range("H65").Value = 3127.8
$ws.Range("I65").Value = 3000.6667
$ws.Range("K65").Value = 15003.3335
$ws.Range("M65").Value = -11883.3335

$ws.Range("H116").Value = 2555.577
$ws.Range("I116").Value = 2009.1333
$ws.Range("K116").Value = 2009.1333
$ws.Range("M116").Value = 1432.8667

$ws.Range("H129").Value = 839.02325
$ws.Range("J129").Value = 1023.84375
$ws.Range("L129").Value = 3071.53125
$ws.Range("N129").Value = -13071.53125

$ws.Range("H132").Value = 7758079.5
$ws.Range("I132").Value = 12351892
$ws.Range("K132").Value = 37055676
$ws.Range("M132").Value = -37053146

$ws.Range("H137").Value = 1093.9016
$ws.Range("I137").Value = 850.65515
$ws.Range("J137").Value = 1314.3438
$ws.Range("K137").Value = 2551.96545
$ws.Range("L137").Value = 3943.0314
$ws.Range("M137").Value = -1.965450000000146
$ws.Range("N137").Value = -9043.0314

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3833.8965
$ws.Range("I32").Value = 3569.0356
$ws.Range("K32").Value = 3569.0356
$ws.Range("M32").Value = -3282.0356

$ws.Range("H45").Value = 1475.1428
$ws.Range("I45").Value = 1478
$ws.Range("J45").Value = 1471.3334
$ws.Range("K45").Value = 1478
$ws.Range("L45").Value = 1471.3334
$ws.Range("M45").Value = -1101
$ws.Range("N45").Value = -2225.3334

$ws.Range("H132").Value = 1500.5106
$ws.Range("I132").Value = 1205.3513
$ws.Range("J132").Value = 2592.6
$ws.Range("K132").Value = 3616.0539
$ws.Range("L132").Value = 7777.799999999999
$ws.Range("M132").Value = -1086.0539
$ws.Range("N132").Value = -12837.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 267
$ws.Range("I7").Value = 325.5
$ws.Range("J7").Value = 150
$ws.Range("K7").Value = 325.5
$ws.Range("L7").Value = 150
$ws.Range("M7").Value = -212.5
$ws.Range("N7").Value = -376

$ws.Range("H80").Value = 1496.2
$ws.Range("I80").Value = 989
$ws.Range("K80").Value = 989
$ws.Range("M80").Value = 9

$ws.Range("H83").Value = 1496.2
$ws.Range("I83").Value = 989
$ws.Range("K83").Value = 4945
$ws.Range("M83").Value = 47

$ws.Range("H99").Value = 71429704
$ws.Range("I99").Value = 83334430
$ws.Range("K99").Value = 83334430
$ws.Range("M99").Value = -83332932

$ws.Range("H107").Value = 1584.1666
$ws.Range("I107").Value = 1519.091
$ws.Range("J107").Value = 2300
$ws.Range("K107").Value = 1519.091
$ws.Range("L107").Value = 2300
$ws.Range("M107").Value = 400.9090000000001
$ws.Range("N107").Value = -6140

$ws.Range("H134").Value = 4083.9268
$ws.Range("I134").Value = 938.34283
$ws.Range("J134").Value = 22433.166
$ws.Range("K134").Value = 2815.02849
$ws.Range("L134").Value = 67299.49800000001
$ws.Range("M134").Value = -280.0284900000001
$ws.Range("N134").Value = -72369.49800000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 1500.5
$ws.Range("I6").Value = 1500.5
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 1500.5
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -1387.5
$ws.Range("N6").ClearContents()

$ws.Range("H16").Value = 83334780
$ws.Range("I16").Value = 111112380
$ws.Range("K16").Value = 111112380
$ws.Range("M16").Value = -111112093

$ws.Range("H31").Value = 2125.4783
$ws.Range("I31").Value = 2085.7273
$ws.Range("J31").Value = 3000
$ws.Range("K31").Value = 2085.7273
$ws.Range("L31").Value = 3000
$ws.Range("M31").Value = -1790.7273
$ws.Range("N31").Value = -3590

$ws.Range("H34").Value = 2125.4783
$ws.Range("I34").Value = 2085.7273
$ws.Range("J34").Value = 3000
$ws.Range("K34").Value = 2085.7273
$ws.Range("L34").Value = 3000
$ws.Range("M34").Value = -1883.7273
$ws.Range("N34").Value = -3404

$ws.Range("H58").Value = 645.84906
$ws.Range("I58").Value = 605.8611
$ws.Range("J58").Value = 730.5294
$ws.Range("K58").Value = 605.8611
$ws.Range("L58").Value = 730.5294
$ws.Range("M58").Value = -402.8611
$ws.Range("N58").Value = -1136.5294

$ws.Range("H105").Value = 1225.3334
$ws.Range("I105").Value = 1205.8
$ws.Range("J105").Value = 1249.75
$ws.Range("K105").Value = 1205.8
$ws.Range("L105").Value = 1249.75
$ws.Range("M105").Value = 541.2
$ws.Range("N105").Value = -4743.75

$ws.Range("H107").Value = 1535
$ws.Range("I107").Value = 2324.5
$ws.Range("K107").Value = 2324.5
$ws.Range("M107").Value = -404.5

$ws.Range("H113").Value = 83334780
$ws.Range("I113").Value = 111112380
$ws.Range("K113").Value = 111112380
$ws.Range("M113").Value = -111110210

$ws.Range("H134").Value = 1071.697
$ws.Range("I134").Value = 1102.6154
$ws.Range("K134").Value = 3307.8462
$ws.Range("M134").Value = -772.8462

$ws.Range("H136").Value = 645.84906
$ws.Range("I136").Value = 605.8611
$ws.Range("J136").Value = 730.5294
$ws.Range("K136").Value = 1817.5833
$ws.Range("L136").Value = 2191.5882
$ws.Range("M136").Value = 732.4167000000002
$ws.Range("N136").Value = -7291.5882

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 7159.2666
$ws.Range("J107").Value = 10540.3
$ws.Range("L107").Value = 31620.9
$ws.Range("N107").Value = -35460.89999999999

$ws.Range("H108").Value = 1304.4546
$ws.Range("I108").Value = 371.8
$ws.Range("J108").Value = 2081.6667
$ws.Range("K108").Value = 1115.4
$ws.Range("L108").Value = 6245.000100000001
$ws.Range("M108").Value = 1764.6
$ws.Range("N108").Value = -12005.0001

$ws.Range("H110").Value = 9276.5
$ws.Range("I110").Value = 8000
$ws.Range("J110").Value = 9418.333000000001
$ws.Range("K110").Value = 24000
$ws.Range("L110").Value = 28254.999
$ws.Range("M110").Value = -19910
$ws.Range("N110").Value = -36434.999

$ws.Range("H117").Value = 918.0769
$ws.Range("I117").Value = 547
$ws.Range("J117").Value = 1150
$ws.Range("K117").Value = 1641
$ws.Range("L117").Value = 3450
$ws.Range("M117").Value = 1801
$ws.Range("N117").Value = -10334

$ws.Range("H120").Value = 8459.200000000001
$ws.Range("I120").Value = 2300
$ws.Range("J120").Value = 9999
$ws.Range("K120").Value = 6900
$ws.Range("L120").Value = 29997
$ws.Range("M120").Value = -2062
$ws.Range("N120").Value = -39673

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2920.3333
$ws.Range("I80").Value = 1801
$ws.Range("J80").Value = 3480
$ws.Range("K80").Value = 1801
$ws.Range("L80").Value = 3480
$ws.Range("M80").Value = -803
$ws.Range("N80").Value = -5476

$ws.Range("H83").Value = 2920.3333
$ws.Range("I83").Value = 1801
$ws.Range("J83").Value = 3480
$ws.Range("K83").Value = 9005
$ws.Range("L83").Value = 17400
$ws.Range("M83").Value = -4013
$ws.Range("N83").Value = -27384

$ws.Range("H113").Value = 1558.8
$ws.Range("J113").Value = 1612
$ws.Range("L113").Value = 1612
$ws.Range("N113").Value = -5952

$ws.Range("H132").Value = 2371.2727
$ws.Range("I132").Value = 1722.7142
$ws.Range("J132").Value = 3506.25
$ws.Range("K132").Value = 5168.142599999999
$ws.Range("L132").Value = 10518.75
$ws.Range("M132").Value = -2638.142599999999
$ws.Range("N132").Value = -15578.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1374.75
$ws.Range("I46").Value = 874.5
$ws.Range("K46").Value = 874.5
$ws.Range("M46").Value = -686.5

$ws.Range("H93").Value = 601.5454999999999
$ws.Range("J93").Value = 666.3333
$ws.Range("L93").Value = 666.3333
$ws.Range("N93").Value = -3162.3333

$ws.Range("H100").Value = 615.8889
$ws.Range("I100").Value = 438.45834
$ws.Range("J100").Value = 970.75
$ws.Range("K100").Value = 438.45834
$ws.Range("L100").Value = 970.75
$ws.Range("M100").Value = 102.54166
$ws.Range("N100").Value = -2052.75

$ws.Range("H107").Value = 3040
$ws.Range("I107").Value = 3040
$ws.Range("K107").Value = 3040
$ws.Range("M107").Value = -1120

$ws.Range("H132").Value = 16881.477
$ws.Range("I132").Value = 1057.5333
$ws.Range("J132").Value = 52485.35
$ws.Range("K132").Value = 3172.5999
$ws.Range("L132").Value = 157456.05
$ws.Range("M132").Value = -642.5999000000002
$ws.Range("N132").Value = -162516.05

$ws.Range("H136").Value = 1411.5834
$ws.Range("I136").Value = 1190.6875
$ws.Range("J136").Value = 1853.375
$ws.Range("K136").Value = 3572.0625
$ws.Range("L136").Value = 5560.125
$ws.Range("M136").Value = -1022.0625
$ws.Range("N136").Value = -10660.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 3000
$ws.Range("J11").Value = 3000
$ws.Range("L11").Value = 3000
$ws.Range("N11").Value = -3284

$ws.Range("H100").Value = 1101.3334
$ws.Range("I100").Value = 1474.6
$ws.Range("K100").Value = 2949.2
$ws.Range("M100").Value = -2408.2

$ws.Range("H107").Value = 427.65
$ws.Range("I107").Value = 397.2353
$ws.Range("J107").Value = 600
$ws.Range("K107").Value = 1191.7059
$ws.Range("L107").Value = 1800
$ws.Range("M107").Value = 728.2941000000001
$ws.Range("N107").Value = -5640

$ws.Range("H132").Value = 1896.8772
$ws.Range("I132").Value = 1910.5532
$ws.Range("J132").Value = 1832.6
$ws.Range("K132").Value = 5731.6596
$ws.Range("L132").Value = 5497.799999999999
$ws.Range("M132").Value = -3201.6596
$ws.Range("N132").Value = -10557.8
